$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 12 new data rows (index 204-215) after the existing last row (205),
# extending the table from A1:B205 to A1:B217. Column A keeps the same direct
# formatting (bold/border/centered style) as the rest of the index column, so
# it is copied from the existing last-row cell before the new value is written.
$styleSource = $ws.Range("A205")

$styleSource.Copy($ws.Range("A206")) | Out-Null
$ws.Range("A206").Value = 204
$ws.Range("B206").Value = 0.673527324343506

$styleSource.Copy($ws.Range("A207")) | Out-Null
$ws.Range("A207").Value = 205
$ws.Range("B207").Value = 0.7871475579069618

$styleSource.Copy($ws.Range("A208")) | Out-Null
$ws.Range("A208").Value = 206
$ws.Range("B208").Value = 0.6518097941802697

$styleSource.Copy($ws.Range("A209")) | Out-Null
$ws.Range("A209").Value = 207
$ws.Range("B209").Value = 0.4854506742370476

$styleSource.Copy($ws.Range("A210")) | Out-Null
$ws.Range("A210").Value = 208
$ws.Range("B210").Value = 0.5266146202980838

$styleSource.Copy($ws.Range("A211")) | Out-Null
$ws.Range("A211").Value = 209
$ws.Range("B211").Value = 0.4992395822771977

$styleSource.Copy($ws.Range("A212")) | Out-Null
$ws.Range("A212").Value = 210
$ws.Range("B212").Value = 0.3804116394606104

$styleSource.Copy($ws.Range("A213")) | Out-Null
$ws.Range("A213").Value = 211
$ws.Range("B213").Value = 0.3400993612491129

$styleSource.Copy($ws.Range("A214")) | Out-Null
$ws.Range("A214").Value = 212
$ws.Range("B214").Value = 0.5862313697657914

$styleSource.Copy($ws.Range("A215")) | Out-Null
$ws.Range("A215").Value = 213
$ws.Range("B215").Value = 0.4513839602555004

$styleSource.Copy($ws.Range("A216")) | Out-Null
$ws.Range("A216").Value = 214
$ws.Range("B216").Value = 0.411639460610362

$styleSource.Copy($ws.Range("A217")) | Out-Null
$ws.Range("A217").Value = 215
$ws.Range("B217").Value = 0.5308729595457772

